# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# cells produced by the report generator with their newly-generated values.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# G2: Latest HO Xliff Generate Date for 5f969f78-...md
$wsOverview.Range("G2").Value = "2017-01-03 08:11:44"

# --- "zh-cn" sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# H2: Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2017-01-03 08:11:31"
# L2: Correspond Handback DateTime
$wsZhCn.Range("L2").Value = "2017-01-03 08:12:03"

# --- "de-de" sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# H2: Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2017-01-03 08:11:44"
# L2: Correspond Handback DateTime
$wsDeDe.Range("L2").Value = "2017-01-03 08:12:16"
